# Apply the edits described by the diff:
#  - Update several Customer_Reviews (column C) values by +/-1
#  - Reset the custom/bestFit width on columns B and D back to the sheet's
#    standard (default) width, leaving column C's custom width untouched
#  - Update the active cell selection to F4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values in column C (Customer_Reviews) ---
$changes = @{
    2  = 221
    6  = 1286
    9  = 1980
    19 = 879
    20 = 811
    21 = 1793
    22 = 824
    23 = 1682
    24 = 563
    26 = 1787
    27 = 2298
    33 = 1602
    35 = 2022
    36 = 1423
    38 = 2245
    39 = 939
    41 = 1669
    44 = 2117
    49 = 1197
}

foreach ($row in $changes.Keys) {
    $ws.Range("C$row").Value = $changes[$row]
}

# --- Reset columns B and D to the sheet's standard (default) width ---
$standardWidth = $ws.StandardWidth
$ws.Columns.Item(2).ColumnWidth = $standardWidth
$ws.Columns.Item(4).ColumnWidth = $standardWidth

# --- Update the active cell selection ---
$ws.Range("F4").Select()
